# Performance evaluation: add "Layout Duration (s)" column between the
# existing "Layout FPS" and "Exploration min FPS" columns, and append a
# couple of scratch values lower on the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at F. Excel shifts the old F/G/H columns (and their
# column-width formatting) one place to the right to G/H/I, and carries the
# left neighbour's cell style onto the new column's header-row cell - this
# matches the diff's shift of "Exploration min/avg FPS" from F/G to G/H and
# the note column from H to I.
$ws.Columns("F:F").Insert()

# New header for the inserted column.
$ws.Range("F4").Value = "Layout Duration (s)"

# New "Layout Duration (s)" sample values for each data row.
$ws.Range("F5").Value = 5
$ws.Range("F6").Value = 7.5
$ws.Range("F7").Value = 9.5
$ws.Range("F8").Value = 10.5
$ws.Range("F9").Value = 15

$ws.Range("F11").Value = 10
$ws.Range("F12").Value = 15
$ws.Range("F13").Value = 27
$ws.Range("F14").Value = 30
$ws.Range("F15").Value = 20

# A couple of extra scratch values further down column F.
$ws.Range("F25").Value = 30
$ws.Range("F26").Value = 20

# Give the new column its own (non bestFit) width, close to the source width.
$ws.Columns("F:F").ColumnWidth = 16.6

# Move the selection to match where the author left the cursor.
$ws.Range("F20").Select() | Out-Null
